$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Update the OT value for the existing row 91 (VARELA 3018) ---
$ws.Range("E91").Value = "ICD30830056"

# --- 2) Insert a brand-new row at position 92, pushing rows 92:98 down to 93:99 ---
$ws.Rows(92).Insert()

# Fill in the data for the newly inserted row 92 (caso 7260 - Vidal 1861)
# Columns A, B, D, E look numeric/date-like; prefix with an apostrophe so
# Excel stores them as text, matching the rest of the column.
$ws.Range("A92").Value = "'7260"
$ws.Range("B92").Value = "'9/10/2025"
$ws.Range("C92").Value = "Vidal 1861"
$ws.Range("D92").Value = "'13"
$ws.Range("E92").Value = "'809642175"
$ws.Range("F92").Value = "AYKO"
$ws.Range("G92").Value = "Pendiente"
$ws.Range("H92").Value = "Picada"
$ws.Range("I92").Value = 1
$ws.Range("J92").Value = "Cambio"
$ws.Range("K92").Value = "Sin equipos"
$ws.Range("L92").Value = "Pasante"
$ws.Range("M92").Value = -58.458298
$ws.Range("N92").Value = -34.566511
$ws.Range("O92").Value = "Colegiales"
$ws.Range("P92").Value = "Capital Norte"

# --- 3) Append a brand-new row 100 (caso 6182 - Los Patos 2702) ---
$ws.Range("A100").Value = "'6182"
$ws.Range("B100").Value = "'9/17/2025"
$ws.Range("C100").Value = "Los Patos 2702"
$ws.Range("D100").Value = "'4"
$ws.Range("E100").Value = "'809818308"
$ws.Range("F100").Value = "AYKO"
$ws.Range("G100").Value = "Pendiente"
$ws.Range("H100").Value = "Sacar PRFV del cantero, colocar en vereda y aplomar"
$ws.Range("I100").Value = 1
$ws.Range("J100").Value = "Cambio"
$ws.Range("K100").Value = "Sin equipos"
$ws.Range("L100").Value = "Pasante"
$ws.Range("M100").Value = -58.399262
$ws.Range("N100").Value = -34.639685
$ws.Range("O100").Value = "San Telmo"
$ws.Range("P100").Value = "Capital Sur"
